$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1644.8889
$ws.Range("I88").Value = 400
$ws.Range("J88").Value = 1800.5
$ws.Range("K88").Value = 400
$ws.Range("L88").Value = 1800.5
$ws.Range("M88").Value = 6
$ws.Range("N88").Value = -2612.5

$ws.Range("H91").Value = 1644.8889
$ws.Range("I91").Value = 400
$ws.Range("J91").Value = 1800.5
$ws.Range("K91").Value = 400
$ws.Range("L91").Value = 1800.5
$ws.Range("M91").Value = 1004
$ws.Range("N91").Value = -4608.5

$ws.Range("H112").Value = 1143.5
$ws.Range("J112").Value = 1158.9474
$ws.Range("L112").Value = 3476.8422
$ws.Range("N112").Value = -5692.8422


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5397.7046
$ws.Range("I32").Value = 4419.024
$ws.Range("K32").Value = 4419.024
$ws.Range("M32").Value = -4132.024


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2886.9
$ws.Range("I105").Value = 2799.875
$ws.Range("J105").Value = 3235
$ws.Range("K105").Value = 2799.875
$ws.Range("L105").Value = 3235
$ws.Range("M105").Value = -1052.875
$ws.Range("N105").Value = -6729

$ws.Range("H113").Value = 1726.6666
$ws.Range("I113").Value = 1726.6666
$ws.Range("K113").Value = 1726.6666
$ws.Range("M113").Value = 443.3334


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2075.8
$ws.Range("I16").Value = 2075.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2075.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1788.8
$ws.Range("N16").Value = ""

$ws.Range("H58").Value = 2610.5217
$ws.Range("I58").Value = 2554.2
$ws.Range("J58").Value = 2653.8462
$ws.Range("K58").Value = 2554.2
$ws.Range("L58").Value = 2653.8462
$ws.Range("M58").Value = -2351.2
$ws.Range("N58").Value = -3059.8462

$ws.Range("H99").Value = 9157.5
$ws.Range("I99").Value = 2067.4443
$ws.Range("J99").Value = 21919.6
$ws.Range("K99").Value = 2067.4443
$ws.Range("L99").Value = 21919.6
$ws.Range("M99").Value = -569.4443000000001
$ws.Range("N99").Value = -24915.6

$ws.Range("H113").Value = 2075.8
$ws.Range("I113").Value = 2075.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2075.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 94.19999999999982
$ws.Range("N113").Value = ""

$ws.Range("H126").Value = 9157.5
$ws.Range("I126").Value = 2067.4443
$ws.Range("J126").Value = 21919.6
$ws.Range("K126").Value = 6202.3329
$ws.Range("L126").Value = 65758.79999999999
$ws.Range("M126").Value = -3732.3329
$ws.Range("N126").Value = -70698.79999999999

$ws.Range("H134").Value = 2765.647
$ws.Range("I134").Value = 2932.5334
$ws.Range("J134").Value = 1514
$ws.Range("K134").Value = 8797.600199999999
$ws.Range("L134").Value = 4542
$ws.Range("M134").Value = -6262.600199999999
$ws.Range("N134").Value = -9612

$ws.Range("H136").Value = 2610.5217
$ws.Range("I136").Value = 2554.2
$ws.Range("J136").Value = 2653.8462
$ws.Range("K136").Value = 7662.599999999999
$ws.Range("L136").Value = 7961.5386
$ws.Range("M136").Value = -5112.599999999999
$ws.Range("N136").Value = -13061.5386

$ws.Range("H140").Value = 75826.22
$ws.Range("J140").Value = 75826.22
$ws.Range("L140").Value = 75826.22
$ws.Range("N140").Value = -86186.22


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 3633.3333
$ws.Range("I120").Value = 2950
$ws.Range("K120").Value = 8850
$ws.Range("M120").Value = -4012

$ws.Range("H122").Value = 1385.3478
$ws.Range("I122").Value = 498.10526
$ws.Range("J122").Value = 5599.75
$ws.Range("K122").Value = 4482.94734
$ws.Range("L122").Value = 50397.75
$ws.Range("M122").Value = -2032.94734
$ws.Range("N122").Value = -55297.75

$ws.Range("H132").Value = 914.8929000000001
$ws.Range("I132").Value = 846.0476
$ws.Range("J132").Value = 1121.4286
$ws.Range("K132").Value = 7614.4284
$ws.Range("L132").Value = 10092.8574
$ws.Range("M132").Value = -5084.4284
$ws.Range("N132").Value = -15152.8574


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 311.0625
$ws.Range("I107").Value = 142.125
$ws.Range("J107").Value = 480
$ws.Range("K107").Value = 142.125
$ws.Range("L107").Value = 480
$ws.Range("M107").Value = 1777.875
$ws.Range("N107").Value = -4320

$ws.Range("H113").Value = 1021.55554
$ws.Range("I113").Value = 1046.125
$ws.Range("K113").Value = 1046.125
$ws.Range("M113").Value = 1123.875

$ws.Range("H138").Value = 46286
$ws.Range("J138").Value = 46286
$ws.Range("L138").Value = 46286
$ws.Range("N138").Value = -56566

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1166.8
$ws.Range("I61").Value = 1354.8
$ws.Range("J61").Value = 978.8
$ws.Range("K61").Value = 1354.8
$ws.Range("L61").Value = 978.8
$ws.Range("M61").Value = -1152.8
$ws.Range("N61").Value = -1382.8

$ws.Range("H113").Value = 1166.8
$ws.Range("I113").Value = 1354.8
$ws.Range("J113").Value = 978.8
$ws.Range("K113").Value = 1354.8
$ws.Range("L113").Value = 978.8
$ws.Range("M113").Value = 815.2
$ws.Range("N113").Value = -5318.8

$ws.Range("H132").Value = 12033.05
$ws.Range("I132").Value = 9636.773999999999
$ws.Range("J132").Value = 20286.889
$ws.Range("K132").Value = 28910.322
$ws.Range("L132").Value = 60860.667
$ws.Range("M132").Value = -26380.322
$ws.Range("N132").Value = -65920.667

$ws.Range("H133").Value = 37794.375
$ws.Range("J133").Value = 37794.375
$ws.Range("L133").Value = 37794.375
$ws.Range("N133").Value = -42854.375


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 431.8
$ws.Range("I113").Value = 328
$ws.Range("J113").Value = 501
$ws.Range("K113").Value = 984
$ws.Range("L113").Value = 1503
$ws.Range("M113").Value = 1186
$ws.Range("N113").Value = -5843

$ws.Range("H132").Value = 272151.88
$ws.Range("I132").Value = 527985.4
$ws.Range("J132").Value = 2105.389
$ws.Range("K132").Value = 1583956.2
$ws.Range("L132").Value = 6316.167
$ws.Range("M132").Value = -1581426.2
$ws.Range("N132").Value = -11376.167

